$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing question text and add the answer amount
$ws.Range("A2").Value = "This should be the only question?"
$ws.Range("B2").Value = 1

# Add a new question row
$ws.Range("A3").Value = "Did I lie last time?"
$ws.Range("B3").Value = 2

# Move the active selection as recorded in the workbook
$ws.Range("I8").Select()
